$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '< -0.0425'
$ws.Range('H2').Value = '< -0.0499'
$ws.Range('I2').Value = '< -0.0425'
$ws.Range('J2').Value = '0.0018 -0.0548'
$ws.Range('K2').Value = '< -0.0548'
$ws.Range('L2').Value = '0.0018 -0.0548'
$ws.Range('E3').Value = '0.0865 -0.0206'
$ws.Range('H3').Value = '< -0.028'
$ws.Range('I3').Value = '0.0857 -0.0206'
$ws.Range('J3').Value = '0.034 -0.0329'
$ws.Range('K3').Value = '< -0.0329'
$ws.Range('L3').Value = '0.034 -0.0329'
$ws.Range('E4').Value = '0.7357 0.0121'
$ws.Range('G4').Value = '0.4531 0.0023'
$ws.Range('H4').Value = '0.818 0.0046'
$ws.Range('I4').Value = '0.9462 0.0121'
$ws.Range('J4').Value = '0.5902 -2e-04'
$ws.Range('K4').Value = '0.5075 -3e-04'
$ws.Range('L4').Value = '0.5902 -2e-04'
$ws.Range('F5').Value = '0.2404 -0.0121'
$ws.Range('G5').Value = '0.1388 -0.0098'
$ws.Range('H5').Value = '0.4144 -0.0074'
$ws.Range('I5').Value = '0.5166 0'
$ws.Range('J5').Value = '0.4032 -0.0123'
$ws.Range('K5').Value = '0.3063 -0.0123'
$ws.Range('L5').Value = '0.4032 -0.0123'
$ws.Range('M5').Value = '0.2404 -0.0121'
$ws.Range('N5').Value = '0.2404 -0.0121'
$ws.Range('O5').Value = '0.2404 -0.0121'
$ws.Range('G6').Value = '0.4531 0.0023'
$ws.Range('H6').Value = '0.818 0.0046'
$ws.Range('I6').Value = '0.9462 0.0121'
$ws.Range('J6').Value = '0.5902 -2e-04'
$ws.Range('K6').Value = '0.5075 -3e-04'
$ws.Range('L6').Value = '0.5902 -2e-04'
$ws.Range('H7').Value = '0.575 0.0023'
$ws.Range('I7').Value = '0.6198 0.0098'
$ws.Range('J7').Value = '0.524 -0.0026'
$ws.Range('K7').Value = '0.5219 -0.0026'
$ws.Range('L7').Value = '0.524 -0.0026'
$ws.Range('M7').Value = '0.4995 -0.0023'
$ws.Range('N7').Value = '0.4995 -0.0023'
$ws.Range('O7').Value = '0.4995 -0.0023'
$ws.Range('I8').Value = '0.303 0.0074'
$ws.Range('J8').Value = '0.3138 -0.0049'
$ws.Range('K8').Value = '0.1709 -0.0049'
$ws.Range('L8').Value = '0.3138 -0.0049'
$ws.Range('M8').Value = '0.1408 -0.0046'
$ws.Range('N8').Value = '0.1408 -0.0046'
$ws.Range('O8').Value = '0.1408 -0.0046'
$ws.Range('J9').Value = '0.2058 -0.0123'
$ws.Range('K9').Value = '0.0307 -0.0123'
$ws.Range('L9').Value = '0.2058 -0.0123'
$ws.Range('M9').Value = '0.0307 -0.0121'
$ws.Range('N9').Value = '0.0307 -0.0121'
$ws.Range('O9').Value = '0.0307 -0.0121'
$ws.Range('K10').Value = '0.4493 0'
$ws.Range('M10').Value = '0.3846 2e-04'
$ws.Range('N10').Value = '0.3846 2e-04'
$ws.Range('O10').Value = '0.3846 2e-04'
$ws.Range('L11').Value = '0.509 0'
$ws.Range('M11').Value = '0.457 3e-04'
$ws.Range('N11').Value = '0.457 3e-04'
$ws.Range('O11').Value = '0.457 3e-04'
$ws.Range('M12').Value = '0.3846 2e-04'
$ws.Range('N12').Value = '0.3846 2e-04'
$ws.Range('O12').Value = '0.3846 2e-04'
